# Update the Coin / Link / Price / Volume(1h) table with the latest
# coinranking.com snapshot (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values such as "0.619" or "232.20" read as plain numbers, so Excel
# would silently convert them (dropping the trailing zero / exact digits).
# Force those specific cells to Text, write the value, then restore the
# "Normal" cell style so no stray number format sticks to the cell.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2: D2='34.915.23', E2='+0.77%'
$ws.Range("D2").Value = '34.915.23'
$ws.Range("E2").Value = '  +0.77%  '

# Row 3: D3='1.843.76', E3='+2.18%'
$ws.Range("D3").Value = '1.843.76'
$ws.Range("E3").Value = '  +2.18%  '

# Row 4: E4='+0.16%'
$ws.Range("E4").Value = '  +0.16%  '

# Row 5: D5='232.20', E5='+0.85%'
Set-TextValue "D5" '232.20'
$ws.Range("E5").Value = '  +0.85%  '

# Row 6: D6='0.619', E6='+3.31%'
Set-TextValue "D6" '0.619'
$ws.Range("E6").Value = '  +3.31%  '

# Row 7: E7='+0.12%'
$ws.Range("E7").Value = '  +0.12%  '

# Row 8: D8='41.05', E8='+6.06%'
Set-TextValue "D8" '41.05'
$ws.Range("E8").Value = '  +6.06%  '

# Row 9: E9='+4.13%'
$ws.Range("E9").Value = '  +4.13%  '

# Row 10: D10='0.0691', E10='+2.14%'
Set-TextValue "D10" '0.0691'
$ws.Range("E10").Value = '  +2.14%  '

# Row 11: D11='0.0983', E11='-0.87%'
Set-TextValue "D11" '0.0983'
$ws.Range("E11").Value = '  -0.87%  '

# Row 12: D12='2.109.29', E12='+2.13%'
$ws.Range("D12").Value = '2.109.29'
$ws.Range("E12").Value = '  +2.13%  '

# Row 13: B13='WrappedEther', C13, D13='1.850.25', E13='+2.31%'
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.850.25'
$ws.Range("E13").Value = '  +2.31%  '

# Row 14: B14='Chainlink', C14, D14='11.37', E14='+5.20%'
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D14" '11.37'
$ws.Range("E14").Value = '  +5.20%  '

# Row 15: E15='+1.99%'
$ws.Range("E15").Value = '  +1.99%  '

# Row 16: E16='+2.88%'
$ws.Range("E16").Value = '  +2.88%  '

# Row 17: D17='34.922.19', E17='+0.84%'
$ws.Range("D17").Value = '34.922.19'
$ws.Range("E17").Value = '  +0.84%  '

# Row 18: D18='69.91', E18='+1.00%'
Set-TextValue "D18" '69.91'
$ws.Range("E18").Value = '  +1.00%  '

# Row 19: D19='0.0₃0789', E19='+1.49%'
$ws.Range("D19").Value = '0.0₃0789'
$ws.Range("E19").Value = '  +1.49%  '

# Row 20: D20='240.30', E20='+0.28%'
Set-TextValue "D20" '240.30'
$ws.Range("E20").Value = '  +0.28%  '

# Row 21: E21='+4.09%'
$ws.Range("E21").Value = '  +4.09%  '

# Row 22: D22='4.75', E22='+2.82%'
Set-TextValue "D22" '4.75'
$ws.Range("E22").Value = '  +2.82%  '

# Row 23: E23='+0.06%'
$ws.Range("E23").Value = '  +0.06%  '

# Row 24: D24='2.26', E24='+1.26%'
Set-TextValue "D24" '2.26'
$ws.Range("E24").Value = '  +1.26%  '

# Row 25: D25='171.99', E25='-0.14%'
Set-TextValue "D25" '171.99'
$ws.Range("E25").Value = '  -0.14%  '

# Row 26: D26='7.86', E26='+2.07%'
Set-TextValue "D26" '7.86'
$ws.Range("E26").Value = '  +2.07%  '

# Row 27: D27='17.42', E27='+2.23%'
Set-TextValue "D27" '17.42'
$ws.Range("E27").Value = '  +2.23%  '

# Row 28: E28='+4.38%'
$ws.Range("E28").Value = '  +4.38%  '

# Row 29: D29='1.63', E29='+8.94%'
Set-TextValue "D29" '1.63'
$ws.Range("E29").Value = '  +8.94%  '

# Row 30: E30='+0.07%'
$ws.Range("E30").Value = '  +0.07%  '

# Row 31: E31='+2.16%'
$ws.Range("E31").Value = '  +2.16%  '

# Row 32: E32='+0.07%'
$ws.Range("E32").Value = '  +0.07%  '

# Row 33: E33='-0.54%'
$ws.Range("E33").Value = '  -0.54%  '

# Row 34: E34='+22.74%'
$ws.Range("E34").Value = '  +22.74%  '

# Row 35: E35='+11.18%'
$ws.Range("E35").Value = '  +11.18%  '

# Row 36: B36='TrustWalletToken', C36, D36='1.24', E36='+2.45%'
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D36" '1.24'
$ws.Range("E36").Value = '  +2.45%  '

# Row 37: B37='ImmutableX', C37, D37='0.743', E37='+8.90%'
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D37" '0.743'
$ws.Range("E37").Value = '  +8.90%  '

# Row 38: E38='+11.79%'
$ws.Range("E38").Value = '  +11.79%  '

# Row 39: D39='89.83', E39='-0.69%'
Set-TextValue "D39" '89.83'
$ws.Range("E39").Value = '  -0.69%  '

# Row 40: D40='1.342.10', E40='+2.65%'
$ws.Range("D40").Value = '1.342.10'
$ws.Range("E40").Value = '  +2.65%  '

# Row 41: E41='+3.27%'
$ws.Range("E41").Value = '  +3.27%  '

# Row 42: D42='14.56', E42='+4.06%'
Set-TextValue "D42" '14.56'
$ws.Range("E42").Value = '  +4.06%  '

# Row 43: B43='HuobiToken', C43, D43='2.41', E43='-1.77%'
$ws.Range("B43").Value = 'HuobiToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D43" '2.41'
$ws.Range("E43").Value = '  -1.77%  '

# Row 44: B44='RenderToken', C44, D44='2.26', E44='+2.83%'
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D44" '2.26'
$ws.Range("E44").Value = '  +2.83%  '

# Row 45: E45='+3.74%'
$ws.Range("E45").Value = '  +3.74%  '

# Row 46: D46='0.0529', E46='+3.41%'
Set-TextValue "D46" '0.0529'
$ws.Range("E46").Value = '  +3.41%  '

# Row 47: D47='6.30', E47='+3.38%'
Set-TextValue "D47" '6.30'
$ws.Range("E47").Value = '  +3.38%  '

# Row 48: D48='2.027.08', E48='+1.80%'
$ws.Range("D48").Value = '2.027.08'
$ws.Range("E48").Value = '  +1.80%  '

# Row 49: E49='+0.15%'
$ws.Range("E49").Value = '  +0.15%  '

# Row 50: D50='3.41', E50='+16.82%'
Set-TextValue "D50" '3.41'
$ws.Range("E50").Value = '  +16.82%  '

# Row 51: D51='0.0669', E51='+0.10%'
Set-TextValue "D51" '0.0669'
$ws.Range("E51").Value = '  +0.10%  '
